# 11_AdminInterface.pptx - "Add files via upload" commit replay.
#
# Two logical changes:
#   1. Every slide's fixed date-and-time footer placeholder is bumped
#      from 2019/5/27 to 2019/5/28.
#   2. On slide 2, the first bullet of the subtitle placeholder is
#      reworded from "In the last section, we discuss Admin Interface."
#      to "This section discusses Admin Interface." (the trailing
#      "Admin Interface." run is left untouched).

$p = $ppt.ActivePresentation

$oldDate = "2019/5/27"
$newDate = "2019/5/28"

# --- 1. Update the date placeholder on every slide -----------------------
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shp = $slide.Shapes.Item($j)
        if ($shp.Type -eq 14 -and $shp.PlaceholderFormat.Type -eq 16) {
            $dtRange = $shp.TextFrame.TextRange
            if ($dtRange.Text -eq $oldDate) {
                $dtRange.Text = $newDate
            }
        }
    }
}

# --- 2. Reword the first bullet on slide 2 --------------------------------
$slide2 = $p.Slides.Item(2)
$bodyShape = $null
for ($j = 1; $j -le $slide2.Shapes.Count; $j++) {
    $candidate = $slide2.Shapes.Item($j)
    if ($candidate.HasTextFrame -and $candidate.TextFrame.HasText) {
        if ($candidate.TextFrame.TextRange.Paragraphs(1, 1).Text -like "In the last section*") {
            $bodyShape = $candidate
        }
    }
}
$para1 = $bodyShape.TextFrame.TextRange.Paragraphs(1, 1)

# "In the last section, we discuss Admin Interface."
#  12345678901234567890123456789012345678901234567890
# Characters 1-12  -> "In the last "
# Characters 13-32 -> "section, we discuss "
# Characters 33-48 -> "Admin Interface." (left untouched)

# Edit the second chunk first so the first chunk's character offsets
# (1-12) stay valid regardless of the length change of the replacement.
$midRun = $para1.Characters(13, 20)
if ($midRun.Text -eq "section, we discuss ") {
    $midRun.Text = "section discusses "
}

$firstRun = $para1.Characters(1, 12)
if ($firstRun.Text -eq "In the last ") {
    $firstRun.Text = "This "
}
